$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($ws, $addr, $text)
    # Force the value to be stored as text, matching the workbook's original
    # inline-string cell type, then restore the default (unstyled) cell style
    # so no stray quote-prefix / number-format style is introduced.
    $ws.Range($addr).Value = "'" + $text
    $ws.Range($addr).Style = "Normal"
}

Set-TextCell $ws "D2" "42.592.40"
Set-TextCell $ws "E2" "  -2.12%  "

Set-TextCell $ws "D3" "2.282.13"
Set-TextCell $ws "E3" "  -3.95%  "

Set-TextCell $ws "E4" "  +0.00%  "

Set-TextCell $ws "D5" "300.29"
Set-TextCell $ws "E5" "  -3.19%  "

Set-TextCell $ws "D6" "97.09"
Set-TextCell $ws "E6" "  -6.80%  "

Set-TextCell $ws "D7" "0.503"
Set-TextCell $ws "E7" "  -1.64%  "

Set-TextCell $ws "E8" "  +0.04%  "

Set-TextCell $ws "D9" "0.497"
Set-TextCell $ws "E9" "  -4.92%  "

Set-TextCell $ws "D10" "33.62"
Set-TextCell $ws "E10" "  -6.05%  "

Set-TextCell $ws "D11" "0.0787"
Set-TextCell $ws "E11" "  -2.54%  "

Set-TextCell $ws "D12" "50.63"
Set-TextCell $ws "E12" "  -4.93%  "

Set-TextCell $ws "D14" "6.63"
Set-TextCell $ws "E14" "  -4.56%  "

Set-TextCell $ws "D15" "2.631.21"
Set-TextCell $ws "E15" "  -4.16%  "

Set-TextCell $ws "D16" "15.20"
Set-TextCell $ws "E16" "  -2.29%  "

Set-TextCell $ws "D17" "2.282.08"
Set-TextCell $ws "E17" "  -3.89%  "

Set-TextCell $ws "D18" "0.786"
Set-TextCell $ws "E18" "  -3.04%  "

Set-TextCell $ws "D19" "42.497.76"
Set-TextCell $ws "E19" "  -2.30%  "

Set-TextCell $ws "D20" "0.0₃0892"
Set-TextCell $ws "E20" "  -2.26%  "

Set-TextCell $ws "D21" "11.40"
Set-TextCell $ws "E21" "  -4.38%  "

Set-TextCell $ws "E22" "  -5.17%  "

Set-TextCell $ws "D23" "66.55"
Set-TextCell $ws "E23" "  -2.54%  "

Set-TextCell $ws "D24" "235.16"
Set-TextCell $ws "E24" "  -2.25%  "

Set-TextCell $ws "E25" "  -5.81%  "

Set-TextCell $ws "D26" "2.48"
Set-TextCell $ws "E26" "  -4.65%  "

Set-TextCell $ws "E27" "  -0.39%  "

Set-TextCell $ws "D28" "24.34"
Set-TextCell $ws "E28" "  -5.77%  "

Set-TextCell $ws "D29" "2.17"
Set-TextCell $ws "E29" "  -6.43%  "

Set-TextCell $ws "D30" "164.70"
Set-TextCell $ws "E30" "  +2.53%  "

Set-TextCell $ws "D31" "33.59"
Set-TextCell $ws "E31" "  -8.54%  "

Set-TextCell $ws "D32" "9.10"
Set-TextCell $ws "E32" "  -4.07%  "

Set-TextCell $ws "D33" "0.999"
Set-TextCell $ws "E33" "  -0.12%  "

Set-TextCell $ws "D34" "4.94"
Set-TextCell $ws "E34" "  -5.71%  "

Set-TextCell $ws "E35" "  -4.61%  "

Set-TextCell $ws "D36" "0.0696"
Set-TextCell $ws "E36" "  -5.26%  "

Set-TextCell $ws "D37" "4.33"
Set-TextCell $ws "E37" "  -7.00%  "

Set-TextCell $ws "D38" "2.82"
Set-TextCell $ws "E38" "  -8.78%  "

Set-TextCell $ws "D39" "16.07"
Set-TextCell $ws "E39" "  -11.86%  "

Set-TextCell $ws "D40" "0.0997"
Set-TextCell $ws "E40" "  -5.63%  "

Set-TextCell $ws "D41" "1.76"
Set-TextCell $ws "E41" "  -8.73%  "

Set-TextCell $ws "E42" "  -3.30%  "

Set-TextCell $ws "E43" "  -9.73%  "

Set-TextCell $ws "D44" "1.954.90"
Set-TextCell $ws "E44" "  -3.86%  "

Set-TextCell $ws "D45" "0.0281"
Set-TextCell $ws "E45" "  -3.18%  "

Set-TextCell $ws "D46" "17.74"
Set-TextCell $ws "E46" "  -9.52%  "

Set-TextCell $ws "D47" "9.66"
Set-TextCell $ws "E47" "  -8.67%  "

Set-TextCell $ws "D48" "2.82"
Set-TextCell $ws "E48" "  -9.49%  "

Set-TextCell $ws "E49" "  -4.61%  "

Set-TextCell $ws "B50" "THORChain"
Set-TextCell $ws "C50" "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextCell $ws "D50" "4.67"
Set-TextCell $ws "E50" "  -1.49%  "

Set-TextCell $ws "B51" "RocketPoolETH"
Set-TextCell $ws "C51" "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextCell $ws "D51" "2.503.22"
Set-TextCell $ws "E51" "  -4.05%  "
